$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells remain text, matching the source data format
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.020.08'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '1.831.01'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").Value = '311.70'
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("D7").Value = '0.4625'
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").Value = '0.3712'
$ws.Range("E8").Value = '  +2.05%  '
$ws.Range("D9").Value = '0.07349'
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("D10").Value = '0.8786'
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").Value = '0.07897'
$ws.Range("E11").Value = '  +4.30%  '
$ws.Range("D12").Value = '19.80'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = '1.851.04'
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("D14").Value = '5.344'
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").Value = '6.546'
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").Value = '91.53'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("D18").Value = '0.000008856'
$ws.Range("E18").Value = '  +2.53%  '
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  -0.61%  '
$ws.Range("D20").Value = '14.81'
$ws.Range("E20").Value = '  +2.26%  '
$ws.Range("D21").Value = '27.037.33'
$ws.Range("E21").Value = '  -1.20%  '
$ws.Range("D22").Value = '5.112'
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("D23").Value = '10.55'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = '1.989.41'
$ws.Range("E24").Value = '  -5.32%  '
$ws.Range("D25").Value = '152.43'
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("D26").Value = '1.845'
$ws.Range("E26").Value = '  -1.12%  '
$ws.Range("D27").Value = '18.47'
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("D28").Value = '2.044'
$ws.Range("E28").Value = '  -2.44%  '
$ws.Range("D29").Value = '5.123'
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").Value = '115.86'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").Value = '0.08903'
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").Value = '2.963'
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").Value = '0.7310'
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").Value = '4.454'
$ws.Range("E34").Value = '  +0.66%  '
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("D36").Value = '2.473'
$ws.Range("E36").Value = '  -1.77%  '
$ws.Range("D37").Value = '1.080'
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("D38").Value = '0.01951'
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("D39").Value = '0.05228'
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("E40").Value = '  +0.92%  '
$ws.Range("D41").Value = '7.119'
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = '0.5179'
$ws.Range("E42").Value = '  -0.54%  '
$ws.Range("D43").Value = '0.1630'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = '8.181'
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("D45").Value = '0.4844'
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("D47").Value = '10.18'
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("D48").Value = '102.35'
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("D49").Value = '1.631'
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("D50").Value = '0.06210'
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").Value = '64.96'
$ws.Range("E51").Value = '  +0.74%  '
